# Rename the first scenario sheet and make it the active/selected sheet again
# (it was previously "시나리오 1. 시계열 분석", the workbook had
# "시나리오 2.2. 시계열 딥러닝(multi_LSTM)" selected).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "시나리오 1. 시계열 (Auto_Arima)"
$ws1.Activate()
